$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values per repull of data / mean calculation
$updates = @{
    "F5"  = -1
    "F6"  = -7
    "F7"  = -1
    "F12" = 1
    "F13" = -4
    "F23" = -1
    "F28" = 2
    "F30" = 1
    "F35" = -1
    "F37" = -4
    "F41" = -2
    "F43" = -2
    "F45" = -1
    "F47" = 2
    "F49" = -2
    "F50" = -3
    "F52" = -2
    "F53" = -4
    "F54" = -3
    "F56" = 7
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
